$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.749.91"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.673.52"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "598.92"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "158.08"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.614"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +4.21%  "
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("E10").Value = "  -0.51%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.83"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("E12").Value = "  -0.01%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "29.15"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.94%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.0000200"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "3.149.05"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "65.610.52"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "2.656.10"
$ws.Range("E17").Value = "  -1.19%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "12.71"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.51"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "352.02"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  +0.01%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "69.22"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.67%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.0000114"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.43%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.74"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("E27").Value = "  -3.27%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.166"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.29%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.02"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("E30").Value = "  -0.09%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "534.10"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -3.28%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.79"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.31%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.50"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("E35").Value = "  -0.04%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.424"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.23%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "20.65"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  -0.07%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "156.84"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.37%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.62%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.00%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "162.95"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("E43").Value = "  -1.01%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.35"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +2.71%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0612"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.74%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "22.75"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.15%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.642"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").Value = "0.0₆0254"
$ws.Range("E49").Value = "  +6.71%  "
$ws.Range("E50").Value = "  -1.10%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "19.93"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.18%  "
